# Fruta / hortaliza, semanal
#
# The underlying data rows (4..28) are re-shuffled: each row's full content
# (dates, volumes, prices, quality, origin, etc.) is relocated to a
# different row position as part of the weekly consolidation. We capture
# every source row verbatim first (so we never read an already-overwritten
# row), then write each row's captured values into its new destination row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of destination row -> source row (both are within the same table,
# rows 4 through 28, columns A through R).
$rowMap = @{
    4 = 21
    5 = 22
    6 = 23
    7 = 24
    8 = 15
    9 = 16
    10 = 12
    11 = 13
    12 = 25
    13 = 26
    14 = 17
    15 = 18
    16 = 10
    17 = 11
    18 = 19
    19 = 20
    20 = 14
    21 = 6
    22 = 7
    23 = 4
    24 = 5
    25 = 27
    26 = 28
    27 = 8
    28 = 9
}

# Step 1: snapshot every source row (A:R) before making any changes.
$snapshot = @{}
for ($r = 4; $r -le 28; $r++) {
    $rng = $ws.Range("A" + $r + ":R" + $r)
    $snapshot[$r] = $rng.Value()
}

# Step 2: write each row's captured values into its destination row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $destRange = $ws.Range("A" + $destRow + ":R" + $destRow)
    $destRange.Value = $snapshot[$srcRow]
}
